# Update scraped Universalis market-price figures (Tonberry) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables. Columns:
#   H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#   K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 18
$ws.Range("H18").Value = 12984.25
$ws.Range("I18").Value = 2000
$ws.Range("K18").Value = 2000
$ws.Range("M18").Value = -1716

# row 98
$ws.Range("H98").Value = 2600.4666
$ws.Range("I98").Value = 2769.7693
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 2769.7693
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = -1271.7693
$ws.Range("N98").Value = -4496

# row 103
$ws.Range("H103").Value = 1118.3125
$ws.Range("I103").Value = 982.8333
$ws.Range("J103").Value = 1524.75
$ws.Range("K103").Value = 2948.4999
$ws.Range("L103").Value = 4574.25
$ws.Range("M103").Value = -2362.4999
$ws.Range("N103").Value = -5746.25

# row 122
$ws.Range("H122").Value = 2600.4666
$ws.Range("I122").Value = 2769.7693
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 8309.3079
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -5859.3079
$ws.Range("N122").Value = -9400

# row 132
$ws.Range("H132").Value = 1019.6316
$ws.Range("I132").Value = 953.02856
$ws.Range("J132").Value = 1796.6666
$ws.Range("K132").Value = 2859.08568
$ws.Range("L132").Value = 5389.9998
$ws.Range("M132").Value = -329.0856800000001
$ws.Range("N132").Value = -10449.9998

# row 141
$ws.Range("H141").Value = 2342.8
$ws.Range("I141").Value = 1095.25
$ws.Range("J141").Value = 7333
$ws.Range("K141").Value = 3285.75
$ws.Range("L141").Value = 21999
$ws.Range("M141").Value = 1894.25
$ws.Range("N141").Value = -32359

$ws = $wb.Worksheets.Item("ARM")
# row 5
$ws.Range("H5").Value = 2199.5
$ws.Range("I5").Value = 174.25
$ws.Range("K5").Value = 174.25
$ws.Range("M5").Value = -62.25

# row 32
$ws.Range("H32").Value = 4397.59
$ws.Range("I32").Value = 4209.4375
$ws.Range("J32").Value = 5257.7144
$ws.Range("K32").Value = 4209.4375
$ws.Range("L32").Value = 5257.7144
$ws.Range("M32").Value = -3922.4375
$ws.Range("N32").Value = -5831.7144

# row 74
$ws.Range("H74").Value = 1673.5883
$ws.Range("I74").Value = 1562.7778
$ws.Range("J74").Value = 1798.25
$ws.Range("K74").Value = 1562.7778
$ws.Range("L74").Value = 1798.25
$ws.Range("M74").Value = -688.7778000000001
$ws.Range("N74").Value = -3546.25

# row 77
$ws.Range("H77").Value = 1673.5883
$ws.Range("I77").Value = 1562.7778
$ws.Range("J77").Value = 1798.25
$ws.Range("K77").Value = 7813.889
$ws.Range("L77").Value = 8991.25
$ws.Range("M77").Value = -3445.889
$ws.Range("N77").Value = -17727.25

# row 102
$ws.Range("H102").Value = 2115
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = $null

# row 132
$ws.Range("H132").Value = 1942.3
$ws.Range("I132").Value = 1259.0435
$ws.Range("K132").Value = 3777.1305
$ws.Range("M132").Value = -1247.1305

$ws = $wb.Worksheets.Item("BSM")
# row 4
$ws.Range("H4").Value = 2199.5
$ws.Range("I4").Value = 174.25
$ws.Range("K4").Value = 174.25
$ws.Range("M4").Value = -59.25

# row 99
$ws.Range("H99").Value = 1954.4706
$ws.Range("I99").Value = 1737.3334
$ws.Range("K99").Value = 1737.3334
$ws.Range("M99").Value = -239.3334

$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 350
$ws.Range("I7").Value = 150
$ws.Range("J7").Value = 450
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 450
$ws.Range("M7").Value = -37
$ws.Range("N7").Value = -676

# row 16
$ws.Range("H16").Value = 1844.1111
$ws.Range("J16").Value = 2749.75
$ws.Range("L16").Value = 2749.75
$ws.Range("N16").Value = -3323.75

# row 22
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null

# row 31
$ws.Range("H31").Value = 2255.0667
$ws.Range("I31").Value = 2363.3
$ws.Range("J31").Value = 2038.6
$ws.Range("K31").Value = 2363.3
$ws.Range("L31").Value = 2038.6
$ws.Range("M31").Value = -2068.3
$ws.Range("N31").Value = -2628.6

# row 34
$ws.Range("H34").Value = 2255.0667
$ws.Range("I34").Value = 2363.3
$ws.Range("J34").Value = 2038.6
$ws.Range("K34").Value = 2363.3
$ws.Range("L34").Value = 2038.6
$ws.Range("M34").Value = -2161.3
$ws.Range("N34").Value = -2442.6

# row 69
$ws.Range("H69").Value = 52856.6
$ws.Range("I69").Value = 11020.5
$ws.Range("K69").Value = 11020.5
$ws.Range("M69").Value = -10271.5

# row 72
$ws.Range("H72").Value = 52856.6
$ws.Range("I72").Value = 11020.5
$ws.Range("K72").Value = 33061.5
$ws.Range("M72").Value = -29317.5

# row 88
$ws.Range("H88").Value = 43333.332
$ws.Range("J88").Value = 43333.332
$ws.Range("L88").Value = 43333.332
$ws.Range("N88").Value = -44145.332

# row 91
$ws.Range("H91").Value = 43333.332
$ws.Range("J91").Value = 43333.332
$ws.Range("L91").Value = 43333.332
$ws.Range("N91").Value = -46141.332

# row 113
$ws.Range("H113").Value = 1844.1111
$ws.Range("J113").Value = 2749.75
$ws.Range("L113").Value = 2749.75
$ws.Range("N113").Value = -7089.75

# row 134
$ws.Range("H134").Value = 2812.9412
$ws.Range("I134").Value = 2531.3845
$ws.Range("K134").Value = 7594.1535
$ws.Range("M134").Value = -5059.1535

$ws = $wb.Worksheets.Item("CUL")
# row 39
$ws.Range("H39").Value = 5333
$ws.Range("J39").Value = 5333
$ws.Range("L39").Value = 15999
$ws.Range("N39").Value = -16587

# row 55
$ws.Range("H55").Value = 12570.2
$ws.Range("J55").Value = 2855.3333
$ws.Range("L55").Value = 8565.999899999999
$ws.Range("N55").Value = -8919.999899999999

# row 131
$ws.Range("H131").Value = 898.6900000000001
$ws.Range("I131").Value = 784.2
$ws.Range("J131").Value = 904.71576
$ws.Range("K131").Value = 2352.6
$ws.Range("L131").Value = 2714.14728
$ws.Range("M131").Value = 2687.4
$ws.Range("N131").Value = -12794.14728

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 2533.3333
$ws.Range("I80").Value = 2550
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2550
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1552
$ws.Range("N80").Value = -4496

# row 83
$ws.Range("H83").Value = 2533.3333
$ws.Range("I83").Value = 2550
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 12750
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -7758
$ws.Range("N83").Value = -22484

# row 132
$ws.Range("H132").Value = 1481876.8
$ws.Range("I132").Value = 2026128.9
$ws.Range("K132").Value = 6078386.699999999
$ws.Range("M132").Value = -6075856.699999999

$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3438.7693
$ws.Range("I7").Value = 2338
$ws.Range("K7").Value = 2338
$ws.Range("M7").Value = -2226

# row 40
$ws.Range("H40").Value = 9094.157999999999
$ws.Range("I40").Value = 12385.1
$ws.Range("K40").Value = 12385.1
$ws.Range("M40").Value = -12249.1

# row 43
$ws.Range("H43").Value = 10808.4
$ws.Range("J43").Value = 10808.4
$ws.Range("L43").Value = 10808.4
$ws.Range("N43").Value = -11194.4

# row 46
$ws.Range("H46").Value = 2114.2856
$ws.Range("I46").Value = 1700
$ws.Range("J46").Value = 2666.6667
$ws.Range("K46").Value = 1700
$ws.Range("L46").Value = 2666.6667
$ws.Range("M46").Value = -1512
$ws.Range("N46").Value = -3042.6667

# row 100
$ws.Range("H100").Value = 1349.125
$ws.Range("I100").Value = 1256.1428
$ws.Range("K100").Value = 1256.1428
$ws.Range("M100").Value = -715.1428000000001

# row 122
$ws.Range("H122").Value = 5399.2856
$ws.Range("I122").Value = 3928.3333
$ws.Range("J122").Value = 6502.5
$ws.Range("K122").Value = 11784.9999
$ws.Range("L122").Value = 19507.5
$ws.Range("M122").Value = -9334.999899999999
$ws.Range("N122").Value = -24407.5

# row 126
$ws.Range("H126").Value = 3438.7693
$ws.Range("I126").Value = 2338
$ws.Range("K126").Value = 7014
$ws.Range("M126").Value = -4544

# row 136
$ws.Range("H136").Value = 2729.8125
$ws.Range("I136").Value = 1486.5834
$ws.Range("K136").Value = 4459.7502
$ws.Range("M136").Value = -1909.7502

$ws = $wb.Worksheets.Item("WVR")
# row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null

# row 126
$ws.Range("H126").Value = 3886.7222
$ws.Range("I126").Value = 2272.4285
$ws.Range("K126").Value = 6817.2855
$ws.Range("M126").Value = -4347.2855
